$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 180; $row -le 234; $row++) {
    $name = "B{0:D3}" -f ($row - 1)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = "T1"
    $ws.Cells.Item($row, 3).Value = "T4"
    $ws.Cells.Item($row, 4).Value = "T12"
    $ws.Cells.Item($row, 5).Value = "T12"
}
